# feat: add 2022-Q3 data
#
# 1. Duplicate the existing "2022-Q2" fund-holdings sheet to create a new
#    "2022-Q3" sheet positioned right before it (so the sheet order becomes
#    总计, 2022-Q3, 2022-Q2), then overwrite its data with the Q3 figures.
# 2. On the "总计" (total) summary sheet, insert the new 2022-Q3 summary row
#    above the existing 2022-Q2 row (shifting it down one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" sheet by copying "2022-Q2", placed before it
# ---------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2.Copy($null, $sheetQ2)

$newCopy = $wb.Worksheets.Item(3)
$newCopy.Name = "2022-Q2-tmp"
$sheetQ2.Name = "2022-Q3"
$newCopy.Name = "2022-Q2"

$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Replicate the data-row style (border/font/quote handling) from the
# template row (row 2, copied from the old Q2 sheet) down through row 10,
# so every new fund row keeps the same look as the existing rows.
for ($r = 3; $r -le 10; $r++) {
    $wsQ3.Range("A2").Copy($wsQ3.Range("A$r"))
}

# Fund holdings for 2022-Q3: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$q3Data = @(
    @(0, "001227", "中邮信息产业灵活配置混合",              "5.91", "85.43", "3.85", "0.2275", 4),
    @(1, "010296", "万家互联互通中国优势量化策略混合A",      "4.22", "86.53", "5.27", "0.2224", 10),
    @(2, "005310", "广发电子信息传媒股票A",                  "1.55", "89.36", "4.85", "0.0752", 1),
    @(3, "001275", "中邮创新优势灵活配置混合",                "1.04", "86.62", "3.45", "0.0359", 8),
    @(4, "010447", "中邮未来成长混合A",                      "0.43", "91.79", "5.70", "0.0245", 1),
    @(5, "010297", "万家互联互通中国优势量化策略混合C",      "0.46", "86.53", "5.27", "0.0242", 10),
    @(6, "010236", "广发电子信息传媒股票C",                  "0.13", "89.36", "4.85", "0.0063", 1),
    @(7, "010404", "博道盛利6个月持有期混合",                "1.10", "41.15", "0.45", "0.0050", 5),
    @(8, "010448", "中邮未来成长混合C",                      "0.06", "91.79", "5.70", "0.0034", 1)
)

# Columns B-G carry text values in this workbook (fund codes keep leading
# zeros, figures are stored as plain text) -- force text formatting before
# assigning so Excel doesn't reinterpret them as numbers.
$wsQ3.Range("B2:G10").NumberFormat = "@"

$r = 2
foreach ($item in $q3Data) {
    $wsQ3.Range("A$r").Value = $item[0]
    $wsQ3.Range("B$r").Value = $item[1]
    $wsQ3.Range("C$r").Value = $item[2]
    $wsQ3.Range("D$r").Value = $item[3]
    $wsQ3.Range("E$r").Value = $item[4]
    $wsQ3.Range("F$r").Value = $item[5]
    $wsQ3.Range("G$r").Value = $item[6]
    $wsQ3.Range("H$r").Value = $item[7]
    $r++
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new Q3 row
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Push the existing 2022-Q2 summary row (currently row 2) down to row 3,
# bumping its index counter, while keeping its original style.
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.13

# Write the new 2022-Q3 summary data into row 2.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 0.62

# Restore the originally active sheet ("总计").
$wsTotal.Activate()
